$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.943.20"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "1.621.48"
$ws.Range("E3").Value = "  -1.03%  "

$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.50"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.04%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.497"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.54%  "

$ws.Range("E7").Value = "  -0.32%  "

$ws.Range("E8").Value = "  -1.03%  "

$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.39"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0793"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.31%  "

$ws.Range("D12").Value = "1.846.54"
$ws.Range("E12").Value = "  -1.08%  "

$ws.Range("E13").Value = "  -1.56%  "

$ws.Range("D14").Value = "1.608.04"
$ws.Range("E14").Value = "  -4.84%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.46%  "

$ws.Range("D16").Value = "25.956.38"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.71"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.98%  "

$ws.Range("E18").Value = "  -1.37%  "

$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "192.14"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.43%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.25"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.69%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.53"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.00%  "

$ws.Range("E23").Value = "  -2.22%  "

$ws.Range("E24").Value = "  +0.68%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.96"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.04%  "

$ws.Range("E26").Value = "  -0.33%  "

$ws.Range("E27").Value = "  -3.55%  "

$ws.Range("E28").Value = "  -2.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.22"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.29%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.23"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0478"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.51%  "

$ws.Range("E32").Value = "  -1.50%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.09"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.70%  "

$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.41"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.69%  "

$ws.Range("B35").Value = "LidoDAOToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.49"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.12%  "

$ws.Range("D36").Value = "1.126.20"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("E37").Value = "  -3.79%  "

$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("E39").Value = "  -2.05%  "

$ws.Range("E40").Value = "  -1.76%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.82"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.21%  "

$ws.Range("D42").Value = "1.758.19"
$ws.Range("E42").Value = "  -0.85%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.756"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.69%  "

$ws.Range("E44").Value = "  -3.88%  "

$ws.Range("E45").Value = "  -0.83%  "

$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.44%  "

$ws.Range("E48").Value = "  -1.92%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.49"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("E51").Value = "  -0.14%  "

